# Apply edit: add a "2022-Q1" fund-holdings detail sheet immediately
# before the existing "总计" (totals) summary sheet, and prepend a
# matching 2022-Q1 summary row to the "总计" sheet.
#
# NOTE: a worksheet handle captured before a call that changes the
# Worksheets collection (e.g. Worksheets.Add) can end up resolving to a
# different sheet afterwards, because sheet indices shift. To stay safe
# we finish every write through the handle returned by Add() first, and
# only re-acquire the "总计" handle (by name) once all the Add()/insert
# calls that could renumber sheets have already happened.

$wb = $excel.ActiveWorkbook

# Helper: style a header/row-index cell the same way the other
# quarterly sheets already style theirs - bold, thin box border,
# centered/top-aligned - so the new rows look consistent.
function Format-IndexCell($cell) {
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# --- 1. Insert a new "2022-Q1" worksheet right before "总计" --------
$totalSheet = $wb.Worksheets.Item("总计")
$q1Sheet = $wb.Worksheets.Add($totalSheet)
$q1Sheet.Name = "2022-Q1"

# --- 2. Header row for the fund-holdings detail sheet ----------------
$q1Headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $q1Headers.Length; $i++) {
    $headerCell = $q1Sheet.Cells.Item(1, $i + 2)
    $headerCell.Value = $q1Headers[$i]
    Format-IndexCell $headerCell
}

# --- 3. Fund holdings detail rows ------------------------------------
# Columns: index, code, name, scale, totalPosition, positionRatio, marketValue(亿元), rank
$q1Data = @(
    @(0, "010286", "海富通成长价值混合A", "19.35", "92.83", "3.42", "0.6618", 10),
    @(1, "001404", "招商移动互联网产业股票", "13.45", "90.96", "4.41", "0.5931", 8),
    @(2, "519005", "海富通股票混合", "19.71", "93.50", "2.53", "0.4987", 10),
    @(3, "004745", "长盛创新驱动灵活配置混合", "13.89", "82.99", "2.66", "0.3695", 10),
    @(4, "519011", "海富通精选混合", "10.25", "76.05", "3.56", "0.3649", 8),
    @(5, "010116", "民生加银新兴产业混合A", "9.34", "79.32", "3.37", "0.3148", 9),
    @(6, "009651", "海富通成长甄选混合A", "7.87", "93.36", "3.03", "0.2385", 10),
    @(7, "008655", "招商科技创新混合A", "4.72", "90.71", "4.41", "0.2082", 8),
    @(8, "002450", "平安睿享文娱灵活配置混合A", "4.57", "82.22", "3.91", "0.1787", 8),
    @(9, "519056", "海富通内需热点混合", "3.90", "93.87", "4.49", "0.1751", 7),
    @(10, "006058", "民生加银新兴成长混合", "4.88", "85.26", "3.44", "0.1679", 8),
    @(11, "006100", "平安优势产业灵活配置混合A", "4.93", "87.56", "3.33", "0.1642", 9),
    @(12, "010126", "平安价值成长混合A", "5.22", "88.45", "3.07", "0.1603", 10),
    @(13, "550009", "信诚中小盘混合", "5.23", "87.68", "2.83", "0.1480", 10),
    @(14, "519015", "海富通精选贰号混合", "3.73", "76.65", "3.61", "0.1347", 8),
    @(15, "009008", "平安科技创新混合A", "4.23", "87.24", "3.04", "0.1286", 10),
    @(16, "501200", "民生加银科技创新 3 年封闭运作灵活配置混合", "3.18", "85.18", "3.47", "0.1103", 10),
    @(17, "008085", "海富通先进制造股票A", "3.57", "93.07", "3.03", "0.1082", 10),
    @(18, "006101", "平安优势产业灵活配置混合C", "3.11", "87.56", "3.33", "0.1036", 9),
    @(19, "009652", "海富通成长甄选混合C", "3.12", "93.36", "3.03", "0.0945", 10),
    @(20, "004845", "南华瑞盈混合A", "2.90", "82.72", "2.67", "0.0774", 8),
    @(21, "008656", "招商科技创新混合C", "1.57", "90.71", "4.41", "0.0692", 8),
    @(22, "002451", "平安睿享文娱灵活配置混合C", "1.65", "82.22", "3.91", "0.0645", 8),
    @(23, "010287", "海富通成长价值混合C", "1.81", "92.83", "3.42", "0.0619", 10),
    @(24, "519139", "海富通沪港深灵活配置混合", "1.32", "94.37", "4.63", "0.0611", 5),
    @(25, "001815", "华泰柏瑞激励动力灵活配置混合A", "3.98", "39.06", "1.30", "0.0517", 9),
    @(26, "010127", "平安价值成长混合C", "1.56", "88.45", "3.07", "0.0479", 10),
    @(27, "008084", "海富通先进制造股票C", "1.50", "93.07", "3.03", "0.0454", 10),
    @(28, "009009", "平安科技创新混合C", "1.23", "87.24", "3.04", "0.0374", 10),
    @(29, "010117", "民生加银新兴产业混合C", "0.93", "79.32", "3.37", "0.0313", 9),
    @(30, "166109", "信达澳银量化先锋混合（LOF）A", "1.06", "89.44", "2.69", "0.0285", 10),
    @(31, "004332", "恒生前海沪港深新兴产业精选混合", "0.52", "80.98", "3.92", "0.0204", 8),
    @(32, "002082", "华泰柏瑞激励动力灵活配置混合C", "1.26", "39.06", "1.30", "0.0164", 9),
    @(33, "004846", "南华瑞盈混合C", "0.10", "82.72", "2.67", "0.0027", 8),
    @(34, "166110", "信达澳银量化先锋混合（LOF）C", "0.09", "89.44", "2.69", "0.0024", 10)
)

$r = 2
foreach ($row in $q1Data) {
    $idxCell = $q1Sheet.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    Format-IndexCell $idxCell
    $q1Sheet.Cells.Item($r, 2).NumberFormat = "@"
    $q1Sheet.Cells.Item($r, 2).Value = $row[1]
    $q1Sheet.Cells.Item($r, 3).Value = $row[2]
    $q1Sheet.Cells.Item($r, 4).NumberFormat = "@"
    $q1Sheet.Cells.Item($r, 4).Value = $row[3]
    $q1Sheet.Cells.Item($r, 5).NumberFormat = "@"
    $q1Sheet.Cells.Item($r, 5).Value = $row[4]
    $q1Sheet.Cells.Item($r, 6).NumberFormat = "@"
    $q1Sheet.Cells.Item($r, 6).Value = $row[5]
    $q1Sheet.Cells.Item($r, 7).NumberFormat = "@"
    $q1Sheet.Cells.Item($r, 7).Value = $row[6]
    $q1Sheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- 4. Re-acquire the "总计" sheet by name (its index shifted from
#        4 to 5 once "2022-Q1" was inserted ahead of it), then prepend
#        a 2022-Q1 summary row -----------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalIdxCell = $totalSheet.Cells.Item(2, 1)
$totalIdxCell.Value = 0
Format-IndexCell $totalIdxCell
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 35
$totalSheet.Cells.Item(2, 4).Value = 5.54

# --- 5. Renumber the index column (col A) of the rows pushed down ----
#        so it stays a contiguous 0-based sequence after the insert ----
for ($r = 3; $r -le 5; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

